$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value = $replacement
    }
}
